$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").Value = "Alterovitz, Ron"
$ws.Range("I3").Select() | Out-Null
